# Samsung USED pricelist update — full row refresh (rows 2-75) on USED_HIGHEST_ALL,
# plus header styling (bold, thin border, centered/top-aligned) on both sheets.

$wb  = $excel.ActiveWorkbook
$wsUsed = $wb.Worksheets.Item("USED_HIGHEST_ALL")
$wsNew  = $wb.Worksheets.Item("NEW_HIGHEST_ALL")

# ---- Full replacement data for USED_HIGHEST_ALL (Samsung section, rows 2-75) ----
$data = @(
    @("Samsung", "Galaxy Z Fold 3 5G", "256GB", 200),
    @("Samsung", "Galaxy Z Fold 3 5G", "512GB", 250),
    @("Samsung", "Galaxy Z Fold 4 5G", "256GB", 350),
    @("Samsung", "Galaxy Z Fold 4 5G", "512GB", 400),
    @("Samsung", "Galaxy Z Fold 4 5G", "1TB", 500),
    @("Samsung", "Galaxy Z Fold 5 5G", "256GB", 550),
    @("Samsung", "Galaxy Z Fold 5 5G", "512GB", 600),
    @("Samsung", "Galaxy Z Fold 5 5G", "1TB", 650),
    @("Samsung", "Galaxy Z Fold 6 5G", "256GB", 770),
    @("Samsung", "Galaxy Z Fold 6 5G", "512GB", 870),
    @("Samsung", "Galaxy Z Fold 6 5G", "1TB", 970),
    @("Samsung", "Galaxy Z Fold 7 5G", "256GB", 1470),
    @("Samsung", "Galaxy Z Fold 7 5G", "512GB", 1520),
    @("Samsung", "Galaxy Z Fold 7 5G", "1TB", 1620),
    @("Samsung", "Galaxy Z Flip 4 5G", "128GB", 100),
    @("Samsung", "Galaxy Z Flip 4 5G", "256GB", 150),
    @("Samsung", "Galaxy Z Flip 4 5G", "512GB", 200),
    @("Samsung", "Galaxy Z Flip 5 5G", "256GB", 250),
    @("Samsung", "Galaxy Z Flip 5 5G", "512GB", 300),
    @("Samsung", "Galaxy Z Flip 6 5G", "256GB", 400),
    @("Samsung", "Galaxy Z Flip 6 5G", "512GB", 450),
    @("Samsung", "Galaxy Z Flip 7 5G", "256GB", 670),
    @("Samsung", "Galaxy Z Flip 7 5G", "512GB", 770),
    @("Samsung", "Galaxy Z Flip 7 FE 5G", "128GB", 470),
    @("Samsung", "Galaxy Z Flip 7 FE 5G", "256GB", 570),
    @("Samsung", "Galaxy S21 5G", "Base", 120),
    @("Samsung", "Galaxy S21+ 5G", "Base", 150),
    @("Samsung", "Galaxy S21 Ultra 5G", "256GB", 200),
    @("Samsung", "Galaxy S21 Ultra 5G", "512GB", 250),
    @("Samsung", "Galaxy S21 FE 5G", "256GB", 100),
    @("Samsung", "Galaxy S22 5G", "128GB", 150),
    @("Samsung", "Galaxy S22 5G", "256GB", 200),
    @("Samsung", "Galaxy S22+ 5G", "128GB", 200),
    @("Samsung", "Galaxy S22+ 5G", "256GB", 250),
    @("Samsung", "Galaxy S22 Ultra 5G", "256GB", 350),
    @("Samsung", "Galaxy S22 Ultra 5G", "512GB", 400),
    @("Samsung", "Galaxy S23 5G", "128GB", 300),
    @("Samsung", "Galaxy S23 5G", "256GB", 350),
    @("Samsung", "Galaxy S23+ 5G", "256GB", 400),
    @("Samsung", "Galaxy S23+ 5G", "512GB", 450),
    @("Samsung", "Galaxy S23 Ultra 5G", "256GB", 500),
    @("Samsung", "Galaxy S23 Ultra 5G", "512GB", 550),
    @("Samsung", "Galaxy S23 Ultra 5G", "1TB", 600),
    @("Samsung", "Galaxy S23 FE 5G", "256GB", 220),
    @("Samsung", "Galaxy S24 5G", "256GB", 500),
    @("Samsung", "Galaxy S24 5G", "512GB", 550),
    @("Samsung", "Galaxy S24+ 5G", "256GB", 570),
    @("Samsung", "Galaxy S24+ 5G", "512GB", 620),
    @("Samsung", "Galaxy S24 Ultra 5G", "256GB", 700),
    @("Samsung", "Galaxy S24 Ultra 5G", "512GB", 750),
    @("Samsung", "Galaxy S24 Ultra 5G", "1TB", 850),
    @("Samsung", "Galaxy S24 FE 5G", "256GB", 370),
    @("Samsung", "Galaxy S24 FE 5G", "512GB", 420),
    @("Samsung", "Galaxy S25 5G", "256GB", 600),
    @("Samsung", "Galaxy S25 5G", "512GB", 700),
    @("Samsung", "Galaxy S25+ 5G", "256GB", 750),
    @("Samsung", "Galaxy S25+ 5G", "512GB", 850),
    @("Samsung", "Galaxy S25 Ultra 5G", "256GB", 850),
    @("Samsung", "Galaxy S25 Ultra 5G", "512GB", 1050),
    @("Samsung", "Galaxy S25 Ultra 5G", "1TB", 1100),
    @("Samsung", "Galaxy S25 Edge 5G", "256GB", 570),
    @("Samsung", "Galaxy S25 Edge 5G", "512GB", 670),
    @("Samsung", "Galaxy S25 FE 5G", "128GB", 370),
    @("Samsung", "Galaxy S25 FE 5G", "256GB", 470),
    @("Samsung", "Galaxy S25 FE 5G", "512GB", 520),
    @("Samsung", "Galaxy A36 5G", "8/256GB", 170),
    @("Samsung", "Galaxy A55 5G", "8/128GB", 170),
    @("Samsung", "Galaxy A55 5G", "8/256GB", 220),
    @("Samsung", "Galaxy A56 5G", "12/256GB", 270),
    @("Samsung", "Galaxy A56 5G", "8/256GB", 360),
    @("Samsung", "Galaxy A73 5G", "8/128GB", 100),
    @("Samsung", "Galaxy A73 5G", "8/256GB", 150),
    @("Samsung", "Galaxy Buds 3", "Base", 50),
    @("Samsung", "Galaxy Buds 3 Pro", "Base", 130)
)

$row = 2
foreach ($rec in $data) {
    $wsUsed.Cells.Item($row, 1).Value = $rec[0]
    $wsUsed.Cells.Item($row, 2).Value = $rec[1]
    $wsUsed.Cells.Item($row, 3).Value = $rec[2]
    $wsUsed.Cells.Item($row, 4).Value = $rec[3]
    $row++
}

# ---- Header formatting: bold font, thin border all around, centered + top aligned ----
$headerUsed = $wsUsed.Range("A1:D1")
$headerUsed.Font.Bold = $true
$headerUsed.Borders.LineStyle = 1
$headerUsed.HorizontalAlignment = -4108
$headerUsed.VerticalAlignment = -4160

$headerNew = $wsNew.Range("A1:D1")
$headerNew.Font.Bold = $true
$headerNew.Borders.LineStyle = 1
$headerNew.HorizontalAlignment = -4108
$headerNew.VerticalAlignment = -4160
